# Auto-applies the diff: append 7 new KHL match rows to Matches_SOG,
# refresh as_of_utc timestamps + derived shot stats on Shots_HA / Shots_Summary,
# and bump Meta_ext as_of_utc / build_version.

$wb = $excel.ActiveWorkbook

# --- Matches_SOG: append rows 357-363 -----------------------------------
$wsMatches = $wb.Worksheets.Item("Matches_SOG")

$newMatchRows = @(
    @{ Row=357; Uid="897850"; DateUtc="2025-12-05T17:00:00"; Home="Трактор"; Away="ЦСКА"; SogHome=25; SogAway=19; Source="khl_text" }
    @{ Row=358; Uid="897851"; DateUtc="2025-12-05T17:00:00"; Home="Салават Юлаев"; Away="Ак Барс"; SogHome=38; SogAway=19; Source="khl_text" }
    @{ Row=359; Uid="897855"; DateUtc="2025-12-05T17:00:00"; Home="Автомобилист"; Away="Северсталь"; SogHome=22; SogAway=36; Source="khl_text" }
    @{ Row=360; Uid="897856"; DateUtc="2025-12-05T17:30:00"; Home="Барыс"; Away="Спартак"; SogHome=36; SogAway=25; Source="khl_text" }
    @{ Row=361; Uid="897854"; DateUtc="2025-12-05T18:00:00"; Home="Лада"; Away="Динамо М"; SogHome=18; SogAway=36; Source="khl_text" }
    @{ Row=362; Uid="897853"; DateUtc="2025-12-05T19:30:00"; Home="Нефтехимик"; Away="Металлург Мг"; SogHome=31; SogAway=43; Source="khl_text" }
    @{ Row=363; Uid="897852"; DateUtc="2025-12-05T19:30:00"; Home="Драконы"; Away="Торпедо"; SogHome=28; SogAway=36; Source="khl_text" }
)

foreach ($row in $newMatchRows) {
    # uid/date_utc/home/away/source are text columns in the source data (even
    # though uid looks numeric) - lead with an apostrophe so Excel keeps them
    # typed as text instead of auto-coercing to a number.
    $wsMatches.Cells.Item($row.Row, 1).Value = "'" + $row.Uid
    $wsMatches.Cells.Item($row.Row, 2).Value = $row.DateUtc
    $wsMatches.Cells.Item($row.Row, 3).Value = $row.Home
    $wsMatches.Cells.Item($row.Row, 4).Value = $row.Away
    $wsMatches.Cells.Item($row.Row, 5).Value = $row.SogHome
    $wsMatches.Cells.Item($row.Row, 6).Value = $row.SogAway
    $wsMatches.Cells.Item($row.Row, 7).Value = $row.Source
}

# --- Shots_HA: refresh as_of_utc + updated shot totals ------------------
$wsShotsHA = $wb.Worksheets.Item("Shots_HA")

$shotsHaUpdates = @(
    @{ Row=2; D="2025-12-05T19:30:00Z" }
    @{ Row=3; D="2025-12-05T19:30:00Z"; E=14; G=402; H=436; I=28.7; J=31.1 }
    @{ Row=4; D="2025-12-05T19:30:00Z" }
    @{ Row=5; D="2025-12-05T19:30:00Z"; F=16; K=529; L=471; M=33.1; N=29.4 }
    @{ Row=6; D="2025-12-05T19:30:00Z" }
    @{ Row=7; D="2025-12-05T19:30:00Z"; E=20; G=644; H=620; I=32.2; J=31 }
    @{ Row=8; D="2025-12-05T19:30:00Z"; F=16; K=454; L=481; M=28.4; N=30.1 }
    @{ Row=9; D="2025-12-05T19:30:00Z" }
    @{ Row=10; D="2025-12-05T19:30:00Z"; E=14; G=408; H=488; I=29.1; J=34.9 }
    @{ Row=11; D="2025-12-05T19:30:00Z"; E=14; G=371; H=503; I=26.5 }
    @{ Row=12; D="2025-12-05T19:30:00Z" }
    @{ Row=13; D="2025-12-05T19:30:00Z"; F=16; K=493; L=472; M=30.8; N=29.5 }
    @{ Row=14; D="2025-12-05T19:30:00Z"; E=20; G=632; H=710; J=35.5 }
    @{ Row=15; D="2025-12-05T19:30:00Z" }
    @{ Row=16; D="2025-12-05T19:30:00Z"; E=13; G=365; H=366; I=28.1; J=28.2 }
    @{ Row=17; D="2025-12-05T19:30:00Z"; F=16; K=529; L=424; M=33.1; N=26.5 }
    @{ Row=18; D="2025-12-05T19:30:00Z" }
    @{ Row=19; D="2025-12-05T19:30:00Z"; F=14; K=448; L=465; M=32; N=33.2 }
    @{ Row=20; D="2025-12-05T19:30:00Z"; F=18; K=597; L=558; M=33.2; N=31 }
    @{ Row=21; D="2025-12-05T19:30:00Z"; E=17; G=618; H=487; I=36.4; J=28.6 }
    @{ Row=22; D="2025-12-05T19:30:00Z" }
    @{ Row=23; D="2025-12-05T19:30:00Z"; F=17; K=433; L=450; M=25.5; N=26.5 }
)

foreach ($u in $shotsHaUpdates) {
    $wsShotsHA.Cells.Item($u.Row, 4).Value = $u.D   # as_of_utc
    if ($u.ContainsKey("E")) { $wsShotsHA.Cells.Item($u.Row, 5).Value = $u.E }
    if ($u.ContainsKey("F")) { $wsShotsHA.Cells.Item($u.Row, 6).Value = $u.F }
    if ($u.ContainsKey("G")) { $wsShotsHA.Cells.Item($u.Row, 7).Value = $u.G }
    if ($u.ContainsKey("H")) { $wsShotsHA.Cells.Item($u.Row, 8).Value = $u.H }
    if ($u.ContainsKey("I")) { $wsShotsHA.Cells.Item($u.Row, 9).Value = $u.I }
    if ($u.ContainsKey("J")) { $wsShotsHA.Cells.Item($u.Row, 10).Value = $u.J }
    if ($u.ContainsKey("K")) { $wsShotsHA.Cells.Item($u.Row, 11).Value = $u.K }
    if ($u.ContainsKey("L")) { $wsShotsHA.Cells.Item($u.Row, 12).Value = $u.L }
    if ($u.ContainsKey("M")) { $wsShotsHA.Cells.Item($u.Row, 13).Value = $u.M }
    if ($u.ContainsKey("N")) { $wsShotsHA.Cells.Item($u.Row, 14).Value = $u.N }
}

# --- Shots_Summary: refresh as_of_utc + updated shot totals -------------
$wsShotsSummary = $wb.Worksheets.Item("Shots_Summary")

$shotsSummaryUpdates = @(
    @{ Row=2; D="2025-12-05T19:30:00Z" }
    @{ Row=3; D="2025-12-05T19:30:00Z"; E=33; F=938; G=1035; H=28.4; I=31.4 }
    @{ Row=4; D="2025-12-05T19:30:00Z" }
    @{ Row=5; D="2025-12-05T19:30:00Z"; E=35; F=1162; G=983; H=33.2; I=28.1 }
    @{ Row=6; D="2025-12-05T19:30:00Z" }
    @{ Row=7; D="2025-12-05T19:30:00Z"; E=34; F=1050; G=1078; H=30.9; I=31.7 }
    @{ Row=8; D="2025-12-05T19:30:00Z"; E=32; F=977; G=907; H=30.5; I=28.3 }
    @{ Row=9; D="2025-12-05T19:30:00Z" }
    @{ Row=10; D="2025-12-05T19:30:00Z"; E=33; F=939; G=1177 }
    @{ Row=11; D="2025-12-05T19:30:00Z"; E=32; F=806; G=1186; H=25.2 }
    @{ Row=12; D="2025-12-05T19:30:00Z" }
    @{ Row=13; D="2025-12-05T19:30:00Z"; E=33; F=1135; G=902; H=34.4; I=27.3 }
    @{ Row=14; D="2025-12-05T19:30:00Z"; E=33; F=1006; G=1185; I=35.9 }
    @{ Row=15; D="2025-12-05T19:30:00Z" }
    @{ Row=16; D="2025-12-05T19:30:00Z"; E=33; F=937; G=984; H=28.4; I=29.8 }
    @{ Row=17; D="2025-12-05T19:30:00Z"; E=34; F=1032; G=869; H=30.4; I=25.6 }
    @{ Row=18; D="2025-12-05T19:30:00Z" }
    @{ Row=19; D="2025-12-05T19:30:00Z"; E=33; F=1089; G=1007; H=33; I=30.5 }
    @{ Row=20; D="2025-12-05T19:30:00Z"; E=35; F=1162; G=1088; H=33.2; I=31.1 }
    @{ Row=21; D="2025-12-05T19:30:00Z"; E=34; F=1203; G=1020; H=35.4; I=30 }
    @{ Row=22; D="2025-12-05T19:30:00Z" }
    @{ Row=23; D="2025-12-05T19:30:00Z"; E=33; F=859; G=885; H=26; I=26.8 }
)

foreach ($u in $shotsSummaryUpdates) {
    $wsShotsSummary.Cells.Item($u.Row, 4).Value = $u.D   # as_of_utc
    if ($u.ContainsKey("E")) { $wsShotsSummary.Cells.Item($u.Row, 5).Value = $u.E }
    if ($u.ContainsKey("F")) { $wsShotsSummary.Cells.Item($u.Row, 6).Value = $u.F }
    if ($u.ContainsKey("G")) { $wsShotsSummary.Cells.Item($u.Row, 7).Value = $u.G }
    if ($u.ContainsKey("H")) { $wsShotsSummary.Cells.Item($u.Row, 8).Value = $u.H }
    if ($u.ContainsKey("I")) { $wsShotsSummary.Cells.Item($u.Row, 9).Value = $u.I }
}

# --- Meta_ext: bump as_of_utc + build_version ----------------------------
$wsMeta = $wb.Worksheets.Item("Meta_ext")
$wsMeta.Cells.Item(2, 2).Value = "2025-12-05T19:30:00Z"
$wsMeta.Cells.Item(2, 4).Value = 36

